$wb = $excel.ActiveWorkbook

# New identifiers / timestamps used for this handback report regeneration
$oldGuid1 = "64506efa-3594-4e51-b0fc-082fc6077f59"
$newGuid1 = "4b3aa5b2-b1dc-44a5-9503-b334888c3791"
$oldGuid2 = "e3454ec3-1e0a-4a6c-ad5f-356ff01fb231"
$newGuid2 = "ffff1059b1f1-dd8a-4248-a543-0e746cb27534"

$newHash = "800a2f64f2bb17eceb0b68967f128ad8955f6073"

$newFileName1 = $newGuid1 + ".md"
$newFileName2 = $newGuid2 + ".md"

$newOverviewDate = "2016-09-05 17:11:53"

$newZhXlf = $newGuid1 + "." + $newHash + ".zh-cn.xlf"
$newZhHandoffDate = "2016-09-05 17:11:48"
$newZhHandbackDate = "2016-09-05 17:12:22"

$newDeXlf = $newGuid1 + "." + $newHash + ".de-de.xlf"
$newDeHandoffDate = $newOverviewDate
$newDeHandbackDate = "2016-09-05 17:12:30"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName1
$wsOverview.Range("B2").Value = "e2e\" + $newFileName1
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = $newFileName2
$wsOverview.Range("B3").Value = "e2e\" + $newFileName2
$wsOverview.Range("G3").Value = $newOverviewDate

$i = 0
foreach ($hl in $wsOverview.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $hl.TextToDisplay = "e2e\" + $newFileName1
    } elseif ($i -eq 2) {
        $hl.TextToDisplay = "e2e\" + $newFileName2
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFileName1
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("I2").Value = $newFileName1
$wsZh.Range("J2").Value = $newZhXlf
$wsZh.Range("K2").Value = $newZhHandbackDate

$wsZh.Range("A3").Value = $newFileName2
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $newZhHandoffDate
$wsZh.Range("I3").Value = $newFileName2
$wsZh.Range("J3").Value = $newZhXlf
$wsZh.Range("K3").Value = $newZhHandbackDate

$i = 0
foreach ($hl in $wsZh.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $hl.TextToDisplay = $newFileName1
    } elseif ($i -eq 2) {
        $hl.TextToDisplay = $newFileName1
    } elseif ($i -eq 3) {
        $hl.TextToDisplay = $newFileName2
    } elseif ($i -eq 4) {
        $hl.TextToDisplay = $newFileName2
    }
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFileName1
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newDeHandoffDate
$wsDe.Range("I2").Value = $newFileName1
$wsDe.Range("J2").Value = $newDeXlf
$wsDe.Range("K2").Value = $newDeHandbackDate

$wsDe.Range("A3").Value = $newFileName2
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $newDeHandoffDate
$wsDe.Range("I3").Value = $newFileName2
$wsDe.Range("J3").Value = $newDeXlf
$wsDe.Range("K3").Value = $newDeHandbackDate

$i = 0
foreach ($hl in $wsDe.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $hl.TextToDisplay = $newFileName1
    } elseif ($i -eq 2) {
        $hl.TextToDisplay = $newFileName1
    } elseif ($i -eq 3) {
        $hl.TextToDisplay = $newFileName2
    } elseif ($i -eq 4) {
        $hl.TextToDisplay = $newFileName2
    }
}
